# "save data done" - add a new "Save" column (H) to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 should look like the other header cells (B1:G1),
# so copy the formatting from G1 ("sum") before setting the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Per-row "era data updated" save flags.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
